$wb = $excel.ActiveWorkbook

# Add a new worksheet "Sheet2" after the last existing sheet (CA_FIL)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Sheet2"

# Header row: FromDate / ToDate
$newSheet.Range("A1").Value = "FromDate"
$newSheet.Range("D1").Value = "ToDate"

# Column labels row, left block then right block
$newSheet.Range("A2").Value = "Day_DD"
$newSheet.Range("B2").Value = "Month_MMM"
$newSheet.Range("C2").Value = "Year_YYYY"
$newSheet.Range("D2").Value = "Day_DD"
$newSheet.Range("E2").Value = "Month_MMM"
$newSheet.Range("F2").Value = "Year_YYYY"

# Data row, left block (May 2020) then right block (June 2019)
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "May"
$newSheet.Range("C3").Value = 2020
$newSheet.Range("D3").Value = 1
$newSheet.Range("E3").Value = "June"
$newSheet.Range("F3").Value = 2019

# Make this the active/visible sheet and select C3 (matches tabSelected + activeTab)
$newSheet.Activate()
[void]$newSheet.Range("C3").Select()
